# Updating test files to match the current format in beta
#
# On the "optimization_parameters" sheet, a new "L_curve" parameter row is
# inserted right after the "Model"/"Sigmoid" row, the "Model" label is
# renamed "production_function", and the old "Deletion" row (which has
# become unused) is removed. The active sheet/tab selection also moves
# from "network_weights" to "optimization_parameters".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("optimization_parameters")

# Insert a new row above row 9 (estimate_params), pushing everything
# below (through the old "Deletion" row) down by one.
$ws.Rows.Item(9).Insert()

# Row 8: rename "Model" -> "production_function" (value stays "Sigmoid").
$ws.Range("A8").Value = "production_function"

# New row 9: L_curve parameter, defaulting to 0, formatted like its
# scientific-notation neighbours (alpha/kk_max/MaxIter/etc.).
$ws.Range("A9").Value = "L_curve"
$ws.Range("B9").Value = 0
$ws.Range("B9").NumberFormat = "0.00E+00"

# The old "Deletion" row (Strain metadata no longer used) has shifted
# down to row 17 after the insertion above - remove it entirely.
$ws.Rows.Item(17).Delete()

# Make "optimization_parameters" the active sheet/tab, with A10 selected.
$ws.Activate()
$ws.Range("A10").Select()
